# Applies the latest scraped cryptocurrency price/volume figures to the
# "cryptos" worksheet (coinranking.com export), refreshing the Price (D)
# and Volume(1h) (E) columns, and fixing the row order for a couple of
# coins whose ranking changed position (Chainlink/WrappedEther and
# Stacks/PEPE and Mantle/FirstDigitalUSD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Values that are purely numeric-looking
# text (e.g. "2.86") are prefixed with a leading apostrophe so Excel keeps
# storing them as text (matching the original "inline string" cells)
# instead of silently re-typing them as numbers.
$updates = [ordered]@{
    "D2" = "66.598.96"
    "E2" = "  +1.04%  "
    "D3" = "3.346.87"
    "E3" = "  +1.12%  "
    "D4" = "'0.999"
    "E4" = "  -0.06%  "
    "D5" = "'587.40"
    "E5" = "  +5.62%  "
    "D6" = "'187.21"
    "E6" = "  -0.11%  "
    "D7" = "'0.999"
    "E7" = "  -0.05%  "
    "D8" = "'0.578"
    "E8" = "  -0.38%  "
    "D9" = "'0.184"
    "E9" = "  +1.88%  "
    "D10" = "'0.586"
    "E10" = "  +0.79%  "
    "D11" = "'47.25"
    "E11" = "  +0.92%  "
    "E12" = "  +2.10%  "
    "D13" = "'651.29"
    "E13" = "  +8.35%  "
    "D14" = "3.885.46"
    "E14" = "  +1.34%  "
    "D15" = "'8.54"
    "E15" = "  -0.99%  "
    "D16" = "66.680.97"
    "E16" = "  +1.21%  "
    "E17" = "  +0.59%  "
    "B18" = "WrappedEther"
    "C18" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D18" = "3.347.53"
    "E18" = "  +1.16%  "
    "B19" = "Chainlink"
    "C19" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "D19" = "'17.93"
    "E19" = "  +0.19%  "
    "D20" = "'11.18"
    "E20" = "  +1.24%  "
    "D21" = "'0.903"
    "E21" = "  +0.42%  "
    "D22" = "'17.74"
    "E22" = "  -4.82%  "
    "D23" = "'5.10"
    "E23" = "  +0.78%  "
    "D24" = "'100.29"
    "E24" = "  -0.43%  "
    "D25" = "'4.01"
    "E25" = "  +1.91%  "
    "E26" = "  +2.70%  "
    "D27" = "'9.73"
    "E27" = "  +2.42%  "
    "D28" = "'32.05"
    "E28" = "  +6.01%  "
    "D29" = "'8.63"
    "E29" = "  -0.50%  "
    "D30" = "'6.99"
    "E30" = "  +4.03%  "
    "D31" = "'611.07"
    "E31" = "  +6.90%  "
    "D32" = "'3.93"
    "E32" = "  +1.89%  "
    "D33" = "'11.15"
    "E33" = "  +1.29%  "
    "D34" = "3.878.23"
    "E34" = "  +4.86%  "
    "E35" = "  +1.86%  "
    "E36" = "  +0.09%  "
    "D37" = "'55.99"
    "E37" = "  -1.52%  "
    "D38" = "'2.79"
    "E38" = "  +5.45%  "
    "D39" = "'0.131"
    "E39" = "  +1.81%  "
    "D40" = "'33.51"
    "E40" = "  -0.64%  "
    "B41" = "PEPE"
    "C41" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D41" = "0.0₃0704"
    "E41" = "  -0.33%  "
    "B42" = "Stacks"
    "C42" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D42" = "'3.22"
    "E42" = "  -0.93%  "
    "D43" = "'0.344"
    "E43" = "  +1.85%  "
    "D44" = "'3.39"
    "E44" = "  +0.61%  "
    "D45" = "'0.0421"
    "E45" = "  +0.84%  "
    "E46" = "  -0.15%  "
    "D47" = "'2.57"
    "E47" = "  +0.26%  "
    "B48" = "FirstDigitalUSD"
    "C48" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D48" = "'1.00"
    "E48" = "  +0.53%  "
    "B49" = "Mantle"
    "C49" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D49" = "'1.39"
    "E49" = "  +11.10%  "
    "D50" = "'2.86"
    "E50" = "  -16.97%  "
    "D51" = "'129.95"
    "E51" = "  +5.50%  "
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
